$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.745.85'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.780.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.72%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '356.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.93%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.585'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.28%  '
$ws.Range("E11").Value = '  +2.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0843'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.214.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.773.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.50%  '
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.686.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0968'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.167'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.05%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0449'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0838'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.20%  '
$ws.Range("E40").Value = '  -4.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.55'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.095.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.944'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.82%  '
